$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E14").Value = "IAuthenticationManager"
$ws.Range("F14").Value = "/nacos/plugin-default-impl/nacos-default-auth-plugin/src/main/java/com/alibaba/nacos/plugin/auth/impl/authenticate/IAuthenticationManager.java"
$ws.Range("E15").Value = "AbstractAuthenticationManager"
$ws.Range("F15").Value = "/nacos/plugin-default-impl/nacos-default-auth-plugin/src/main/java/com/alibaba/nacos/plugin/auth/impl/authenticate/AbstractAuthenticationManager.java"

$ws.Range("E15:F15").RowHeight = 102

$ws.Range("E13").Copy()
$ws.Range("E14").PasteSpecial(-4122)
$ws.Range("E15").PasteSpecial(-4122)

$ws.Range("F13").Copy()
$ws.Range("F14").PasteSpecial(-4122)
$ws.Range("F15").PasteSpecial(-4122)

$ws.Range("G15").Select()
